# Apply updated TPM-based NATMI ligand-receptor scores for Lgi1-Rtn4r
# New rows add the "ECs" sending/target cluster in addition to existing MuSCs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgi1"
$ws.Range("C2").Value = "Rtn4r"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.000409
$ws.Range("H2").Value = 0.001227
$ws.Range("I2").Value = 0.164521319388576
$ws.Range("J2").Value = 0.228024530756365
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2214103333333333
$ws.Range("N2").Value = 0.664231
$ws.Range("O2").Value = 0.5896903116237344
$ws.Range("P2").Value = 0.6831206522767569
$ws.Range("Q2").Value = 0.00009055682633333335
$ws.Range("R2").Value = 0.000815011437
$ws.Range("S2").Value = 0.09701662809899735
$ws.Range("T2").Value = 0.1557682661853895

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgi1"
$ws.Range("C3").Value = "Rtn4r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.000409
$ws.Range("H3").Value = 0.001227
$ws.Range("I3").Value = 0.164521319388576
$ws.Range("J3").Value = 0.228024530756365
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.1540585
$ws.Range("N3").Value = 0.308117
$ws.Range("O3").Value = 0.4103096883762655
$ws.Range("P3").Value = 0.3168793477232431
$ws.Range("Q3").Value = 0.0000630099265
$ws.Range("R3").Value = 0.000378059559
$ws.Range("S3").Value = 0.06750469128957869
$ws.Range("T3").Value = 0.07225626457097553

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Lgi1"
$ws.Range("C4").Value = "Rtn4r"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.002077
$ws.Range("H4").Value = 0.004154
$ws.Range("I4").Value = 0.8354786806114239
$ws.Range("J4").Value = 0.771975469243635
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.2214103333333333
$ws.Range("N4").Value = 0.664231
$ws.Range("O4").Value = 0.5896903116237344
$ws.Range("P4").Value = 0.6831206522767569
$ws.Range("Q4").Value = 0.0004598692623333333
$ws.Range("R4").Value = 0.002759215574
$ws.Range("S4").Value = 0.4926736835247371
$ws.Range("T4").Value = 0.5273523860913675

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Lgi1"
$ws.Range("C5").Value = "Rtn4r"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.002077
$ws.Range("H5").Value = 0.004154
$ws.Range("I5").Value = 0.8354786806114239
$ws.Range("J5").Value = 0.771975469243635
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.1540585
$ws.Range("N5").Value = 0.308117
$ws.Range("O5").Value = 0.4103096883762655
$ws.Range("P5").Value = 0.3168793477232431
$ws.Range("Q5").Value = 0.0003199795045
$ws.Range("R5").Value = 0.001279918018
$ws.Range("S5").Value = 0.3428049970866868
$ws.Range("T5").Value = 0.2446230831522676

